# Updates the Extent Report workbook to reflect a newly-run scenario
# (TC_UI_Zlaata_PLP_01 / "This is Product Listing Page feature") in place of
# the previous scenario (TC_UI_Zlaata_PDP_02 / "Product Details Page Feature"),
# along with refreshed run timestamps / durations.

$wb = $excel.ActiveWorkbook

$newScenario = 'TC_UI_Zlaata_PLP_01 |Verify that the "Home" text link on the Product Listing page is clickable.|"TD_UI_Zlaata_PLP_01"'
$newFeature  = 'This is Product Listing Page feature'
$newTag      = '@TC_UI_Zlaata_PLP_01'

# ---- Scenarios sheet ----
$wsScenarios = $wb.Worksheets.Item("Scenarios")
$wsScenarios.Range("B22").Value = $newScenario
$wsScenarios.Range("D22").Value = "14.353 s"
$wsScenarios.Range("E22").Value = $newFeature

# ---- Tags sheet ----
$wsTags = $wb.Worksheets.Item("Tags")
$wsTags.Range("B22").Value = $newTag
$wsTags.Range("B28").Value = $newTag
$wsTags.Range("C28").Value = $newFeature
$wsTags.Range("H28").Value = $newScenario
$wsTags.Range("C29").Value = $newFeature
$wsTags.Range("H29").Value = $newScenario

# ---- Features sheet ----
$wsFeatures = $wb.Worksheets.Item("Features")
$wsFeatures.Range("B22").Value = $newFeature
$wsFeatures.Range("D22").Value = "14.358 s"

# ---- DB Data sheet (hidden helper sheet backing the dashboard charts) ----
$wsDb = $wb.Worksheets.Item("DB Data")
$wsDb.Range("B3").Value = "Oct 10, 2025 1:24:50 pm"
$wsDb.Range("B4").Value = "Oct 10, 2025 1:24:34 pm"
$wsDb.Range("B5").Value = "Oct 10, 2025 1:24:49 pm"
$wsDb.Range("B6").Value = "15.003 s"

# ---- Dashboard sheet protection password was rotated in the source edit.
# The stored value is a one-way hash of the (unknown) plaintext password, so
# we can't reproduce the exact hash bytes; re-protecting keeps the sheet's
# protected state intact without altering any cell data.
$wsDashboard = $wb.Worksheets.Item("Dashboard")
try {
    $wsDashboard.Unprotect("B5F9")
} catch {
}
try {
    $wsDashboard.Protect("FC0D", $true, $true, $true)
} catch {
}

$wb.Save()
